$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet tab from 2020-12-18 to 2020-12-17
$ws.Name = '2020-12-17'

# Replace rows 2-17 (A:C) with the new 2020-12-17 data set
$data = @(
    ,@('2020-12-17', '分散采购关于安吉县2021-2023年度土地及房屋建筑物业务评估服务政府采购项目的成交结果公告', 'http://www.ajztb.com/jyxx/003002/003002008/003002008005/20201217/757b8c0f-e40f-4a17-8003-95736e6bcab9.html')
    ,@('2020-12-17', '自行采购2020年西南片区农村公路提升工程（深王公路、双桥至岭西）挖方路段土、石料地质调查报告编制服务采购项目竞争性磋商', 'http://www.ajztb.com/jyxx/003002/003002009/003002009001/20201217/692213c9-4053-4099-b746-23bfb3c366b3.html')
    ,@('2020-12-17', '自行采购关于安吉县智慧城市二期建设项目——安吉全域智慧旅游建设（一期）采购项目中标结果公告', 'http://www.ajztb.com/jyxx/003002/003002009/003002009005/20201217/eed0b312-740a-4c33-a083-063ad00fe62d.html')
    ,@('2020-12-17', '自行采购关于安吉县经五路（胜利东路-芜园东路）等道路交通设施采购、安装施工工程的更正公告', 'http://www.ajztb.com/jyxx/003002/003002009/003002009003/20201217/569e1318-7784-4fd2-8392-f19a6f44079f.html')
    ,@('2020-12-17', '自行采购关于安吉县经五路（胜利东路-芜园东路）等道路交通设施采购、安装施工工程的公开招标公告', 'http://www.ajztb.com/jyxx/003002/003002009/003002009001/20201217/42cffb83-9bc6-4d31-80e9-e1dd7b28f04e.html')
    ,@('2020-12-17', '自行采购安吉中诚招标代理有限公司关于安吉县交通投资发展集团有限公司“十四五”规划编制服务采购项目的磋商公告（第二次）', 'http://www.ajztb.com/jyxx/003002/003002009/003002009001/20201217/6fdbf2c2-c838-46a6-91b2-cf4eec51fc98.html')
    ,@('2020-12-17', '自行采购关于安吉县现代农业公共服务中心装修设计施工一体化采购项目的中标结果公告', 'http://www.ajztb.com/jyxx/003002/003002009/003002009005/20201217/e5910c20-ec37-4047-8f00-cdf90b8ec2dd.html')
    ,@('2020-12-17', '自行采购关于安吉第一国际城消防设施改造工程项目的成交结果公告', 'http://www.ajztb.com/jyxx/003002/003002009/003002009005/20201217/a683c7c6-9c91-4c2c-a472-b2ecc86c1048.html')
    ,@('2020-12-17', '分散采购关于安吉县2021-2023年度土地及房屋建筑物业务评估服务政府采购项目的成交结果公告', 'http://www.ajztb.com/jyxx/003002/003002008/003002008005/20201217/757b8c0f-e40f-4a17-8003-95736e6bcab9.html')
    ,@('2020-12-17', '自行采购2020年西南片区农村公路提升工程（深王公路、双桥至岭西）挖方路段土、石料地质调查报告编制服务采购项目竞争性磋商', 'http://www.ajztb.com/jyxx/003002/003002009/003002009001/20201217/692213c9-4053-4099-b746-23bfb3c366b3.html')
    ,@('2020-12-17', '自行采购关于安吉县智慧城市二期建设项目——安吉全域智慧旅游建设（一期）采购项目中标结果公告', 'http://www.ajztb.com/jyxx/003002/003002009/003002009005/20201217/eed0b312-740a-4c33-a083-063ad00fe62d.html')
    ,@('2020-12-17', '自行采购关于安吉县经五路（胜利东路-芜园东路）等道路交通设施采购、安装施工工程的更正公告', 'http://www.ajztb.com/jyxx/003002/003002009/003002009003/20201217/569e1318-7784-4fd2-8392-f19a6f44079f.html')
    ,@('2020-12-17', '自行采购关于安吉县经五路（胜利东路-芜园东路）等道路交通设施采购、安装施工工程的公开招标公告', 'http://www.ajztb.com/jyxx/003002/003002009/003002009001/20201217/42cffb83-9bc6-4d31-80e9-e1dd7b28f04e.html')
    ,@('2020-12-17', '自行采购安吉中诚招标代理有限公司关于安吉县交通投资发展集团有限公司“十四五”规划编制服务采购项目的磋商公告（第二次）', 'http://www.ajztb.com/jyxx/003002/003002009/003002009001/20201217/6fdbf2c2-c838-46a6-91b2-cf4eec51fc98.html')
    ,@('2020-12-17', '自行采购关于安吉县现代农业公共服务中心装修设计施工一体化采购项目的中标结果公告', 'http://www.ajztb.com/jyxx/003002/003002009/003002009005/20201217/e5910c20-ec37-4047-8f00-cdf90b8ec2dd.html')
    ,@('2020-12-17', '自行采购关于安吉第一国际城消防设施改造工程项目的成交结果公告', 'http://www.ajztb.com/jyxx/003002/003002009/003002009005/20201217/a683c7c6-9c91-4c2c-a472-b2ecc86c1048.html')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $cellA = $ws.Cells.Item($r, 1)
    # Force text format first so the date-like "2020-12-17" string is stored
    # verbatim instead of being auto-converted into a date serial number,
    # then strip the format back off so no stray style sticks to the cell.
    $cellA.NumberFormat = "@"
    $cellA.Value = $data[$i][0]
    $cellA.ClearFormats()
    $ws.Cells.Item($r, 2).Value = $data[$i][1]
    $ws.Cells.Item($r, 3).Value = $data[$i][2]
}

# Remove the now-stale rows 18-31 that existed in the old (31-row) sheet
$ws.Range("A18:C31").Clear()

$wb.Save()
